$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that wraps nothing at the very
#    start of the document (Title paragraph).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) Re-insert a "_GoBack" bookmark around the single letter "y" in
#    "Using the segments you have written above" (splits that run into
#    three runs: "Using the segments ", "y", "ou have written above").
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("you have written above", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$yRange = $d.Range($r.Start, $r.Start + 1)
$d.Bookmarks.Add("_GoBack", $yRange)

# ---------------------------------------------------------------------
# 3) Add a new list item after "Reads the program all at once into
#    memory" about converting the math operations to functions.
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Reads the program all at once into memory", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertAt = $r2.End
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$newItem = $d.Range($insertAt + 1, $insertAt + 1)
$newItem.Text = "Convert the 4 basic mathematical operations (add, subtract, multiply, divide) to functions and call them from the cases. "
